$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timesheet rows (10-15), matching the style (date format) of the
# existing date column by copying formats from the last populated row.
$xlPasteFormats = -4122

$newRows = @(
    @{ Row = 10; Date = 43863; Hours = 2.4; Note = "Tried updating software versions" },
    @{ Row = 11; Date = 43864; Hours = 2.5; Note = "Head banging" },
    @{ Row = 12; Date = 43865; Hours = 1.8; Note = "Clean up, implement Z+4" },
    @{ Row = 13; Date = 43865; Hours = 3;   Note = "Implement Zip+4, link to Ebay" },
    @{ Row = 14; Date = 43866; Hours = 3;   Note = "More of same" },
    @{ Row = 15; Date = 43875; Hours = 1.4; Note = "Install update on pythoneverywhere" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy the date cell's number format (s="1") down from row 9 so we
    # don't spawn a duplicate style entry.
    $ws.Range("A9").Copy()
    $ws.Range("A$row").PasteSpecial($xlPasteFormats)

    $ws.Range("A$row").Value = $r.Date
    $ws.Range("B$row").Value = $r.Hours
    $ws.Range("C$row").Value = $r.Note
}

$ws.Range("A16").Select()
